$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.294.32"
$ws.Range("E2").Value = "  -3.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.857.08"
$ws.Range("E3").Value = "  -2.82%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.46"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.00"
$ws.Range("E6").Value = "  -2.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.855.32"
$ws.Range("E7").Value = "  -2.88%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -4.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.44"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -3.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  +2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.07"
$ws.Range("E14").Value = "  -4.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.506.28"
$ws.Range("E15").Value = "  -2.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.860.47"
$ws.Range("E16").Value = "  -2.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.396.18"
$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.18"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -4.23%  "

$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.34"
$ws.Range("E22").Value = "  -6.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.736"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("E24").Value = "  -5.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.05"
$ws.Range("E25").Value = "  -3.58%  "

$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.12"
$ws.Range("E27").Value = "  -3.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  -2.32%  "

$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("E30").Value = "  -2.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.009.63"
$ws.Range("E31").Value = "  -2.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.61"
$ws.Range("E32").Value = "  -4.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  -5.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.23"
$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.58"
$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.821.32"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.64"
$ws.Range("E38").Value = "  +8.68%  "

$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("E41").Value = "  -4.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("E43").Value = "  -4.36%  "

$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000301"
$ws.Range("E44").Value = "  +7.95%  "

$ws.Range("E45").Value = "  -5.82%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "421.94"
$ws.Range("E46").Value = "  -4.62%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.06"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.37"
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.41"
$ws.Range("E51").Value = "  -1.56%  "

